# Refresh the cryptos list (prices / 1h volume changes, and a handful of
# re-ranked rows around #41-44) as produced by the scheduled GitHub
# Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.408.24'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '2.642.79'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'600.19"
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').Value = "'154.46"
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').Value = '2.641.85'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').Value = "'0.138"
$ws.Range('E10').Value = '  +5.71%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').Value = "'0.348"
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').Value = "'27.82"
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('D16').Value = '3.123.06'
$ws.Range('E16').Value = '  +1.47%  '
$ws.Range('D17').Value = '68.309.53'
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').Value = '2.643.34'
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').Value = "'11.40"
$ws.Range('E19').Value = '  +3.28%  '
$ws.Range('D20').Value = "'365.63"
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').Value = "'7.39"
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('E24').Value = '  +2.35%  '
$ws.Range('D25').Value = "'73.27"
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = "'9.99"
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D29').Value = "'0.0000104"
$ws.Range('E29').Value = '  +5.68%  '
$ws.Range('D30').Value = "'0.999"
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').Value = "'574.23"
$ws.Range('E31').Value = '  -1.71%  '
$ws.Range('E32').Value = '  +5.11%  '
$ws.Range('E33').Value = '  +4.56%  '
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('E35').Value = '  +3.29%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +3.50%  '
$ws.Range('D38').Value = "'160.04"
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('E39').Value = '  +4.11%  '
$ws.Range('D40').Value = "'19.20"
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = "'0.367"
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = "'5.38"
$ws.Range('E42').Value = '  +3.05%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').Value = "'17.74"
$ws.Range('E43').Value = '  +3.81%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = "'2.62"
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('E45').Value = '  +12.60%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').Value = "'156.88"
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').Value = "'21.89"
$ws.Range('E50').Value = '  +2.30%  '
$ws.Range('E51').Value = '  +1.23%  '

# The apostrophe-prefixed assignments above mark the cells with the
# "number stored as text" quote-prefix style. Clear that formatting
# (cell by cell - ClearFormats on a multi-area Range only affects the
# first area) so the cells end up with no style, matching the source.
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D14').ClearFormats()
$ws.Range('D19').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D30').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D40').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('D43').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D50').ClearFormats()
